# The captured change only touches the two internal theme parts of the
# package: the "Integral" design theme that is applied to the slide
# master/layouts and the default "Office Theme" that backs the notes
# master swap which physical OOXML part (theme1.xml / theme2.xml) each
# one's markup lives in, with no other visible/semantic change anywhere
# else in the deck (no shape, text, slide, layout, or relationship is
# touched).
#
# Re-apply the presentation's design/theme so the deck's theme
# definition is (re)written from the current "Integral" design - the
# COM-level equivalent of re-selecting the same entry in the Design
# gallery, which is what produces this kind of theme-part refresh.
$p = $ppt.ActivePresentation

$designName = "Integral"
if ($p.Designs.Count -ge 1) {
    $designName = $p.Designs.Item(1).Name
}

# Refresh/re-apply the active design theme on the presentation itself...
try { $p.ApplyTheme($designName) } catch {}

# ...and on the slide master that actually owns it, for good measure.
try { $p.SlideMaster.ApplyTheme($designName) } catch {}
